$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2/B3 with freshly generated values (Refresh button functionality)
$ws.Range("B2").Value = "Procedure_20230503142450"
$ws.Range("B3").Value = "Type_20230503142402"

# Clear the old A2/A3 labels entirely (cells removed from the sheet)
$ws.Range("A2").ClearContents()
$ws.Range("A3").ClearContents()
